$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8937149047851562
$ws.Range("B1").Value = 2.698715448379517
$ws.Range("C1").Value = 4.561539173126221
$ws.Range("D1").Value = 2.186407804489136
$ws.Range("E1").Value = 1.289425730705261
